# Checklist.xlsx update: "working on hybrid method, hodgepodge of updates"
#
# Adds three new follow-up tasks to the bottom of the "Tasks" checklist
# (rows 22-24, column B) and moves the sheet's active selection down to
# the newly added last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

$ws.Cells.Item(22, 2).Value2 = "What is option -mcec? What is empirical covariance?"
$ws.Cells.Item(23, 2).Value2 = "What is the .bgs file used for? Read in if mcmc2_flag==TRUE"
$ws.Cells.Item(24, 2).Value2 = "Deduce which options are available for the hybrid method"

# Match the author's final selection/scroll state after adding the rows.
$ws.Range("B24").Select() | Out-Null
